# Remove the first 4 data rows (Cutoff levels 1-4) from each worksheet,
# shifting the remaining data rows up. Header row (row 1) is preserved.
# After the shift, re-number the "Cutoff" index column (A) so it again
# starts at 0 and increases sequentially for the remaining rows.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Rows("2:5").Delete()

    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 2
    }
}
